$wb = $excel.ActiveWorkbook

# ---- Sheet1 ("Sheet1") ----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Cells.Item(1, 4).Value = "initials"

$ws1.Cells.Item(2, 3).Value = "no"
$ws1.Cells.Item(2, 4).Value = "NU"
$ws1.Cells.Item(3, 3).Value = "no"
$ws1.Cells.Item(3, 4).Value = "NU"
$ws1.Cells.Item(4, 3).Value = "no"
$ws1.Cells.Item(4, 4).Value = "NU"

$ws1.Cells.Item(5, 4).Value = "HB"
$ws1.Cells.Item(6, 4).Value = "HB"
$ws1.Cells.Item(7, 4).Value = "HB"

$ws1.Cells.Item(8, 4).Value = "AN"
$ws1.Cells.Item(9, 4).Value = "AN"
$ws1.Cells.Item(10, 4).Value = "AN"

$ws1.Cells.Item(11, 4).Value = "LP"
$ws1.Cells.Item(12, 4).Value = "LP"
$ws1.Cells.Item(13, 4).Value = "LP"

$ws1.Cells.Item(14, 4).Value = "YC"
$ws1.Cells.Item(15, 4).Value = "YC"
$ws1.Cells.Item(16, 4).Value = "YC"

$ws1.Range("C4").Select()
$ws1.Activate()

# ---- Sheet2 ("DeviceSetupLogins") ----
$ws2 = $wb.Worksheets.Item("DeviceSetupLogins")

$ws2.Cells.Item(1, 4).Value = "initials"

$ws2.Cells.Item(2, 3).Value = "no"
$ws2.Cells.Item(2, 4).Value = "AT"
$ws2.Cells.Item(3, 3).Value = "no"
$ws2.Cells.Item(3, 4).Value = "AT"
$ws2.Cells.Item(4, 3).Value = "no"
$ws2.Cells.Item(4, 4).Value = "AT"

$ws2.Cells.Item(5, 4).Value = "NU"
$ws2.Cells.Item(6, 4).Value = "NU"
$ws2.Cells.Item(7, 4).Value = "NU"

$ws2.Cells.Item(8, 4).Value = "HB"
$ws2.Cells.Item(9, 4).Value = "HB"
$ws2.Cells.Item(10, 4).Value = "HB"

$ws2.Cells.Item(11, 4).Value = "AN"
$ws2.Cells.Item(12, 4).Value = "AN"
$ws2.Cells.Item(13, 4).Value = "AN"

$ws2.Cells.Item(14, 4).Value = "LP"
$ws2.Cells.Item(15, 4).Value = "LP"
$ws2.Cells.Item(16, 4).Value = "LP"

$ws2.Cells.Item(17, 4).Value = "YC"
$ws2.Cells.Item(18, 4).Value = "YC"
$ws2.Cells.Item(19, 4).Value = "YC"

$ws2.Range("C4").Select()

# Sheet1 should be the active/tab-selected sheet (not DeviceSetupLogins)
$ws1.Activate()
